$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9473896622657776
$ws.Range("B1").Value = 2.104326725006104
$ws.Range("C1").Value = 4.634419918060303
$ws.Range("D1").Value = 2.965532541275024
$ws.Range("E1").Value = 1.421232223510742
